$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (German)
$ws.Range("D2").Value = 18459
$ws.Range("G2").Value = 28.95454791700525
$ws.Range("L2").Value = 2922
$ws.Range("M2").Value = 0.5467077788920655
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0

# Row 3 (Spanish)
$ws.Range("D3").Value = 1721
$ws.Range("G3").Value = 39.78675188843695
$ws.Range("L3").Value = 3036
$ws.Range("M3").Value = 4.433864442919107
$ws.Range("R3").Value = 13
$ws.Range("S3").Value = 0.01898558555927154

# Row 7 (Chinese)
$ws.Range("D7").Value = 500
$ws.Range("G7").Value = 36.376
$ws.Range("L7").Value = 1189
$ws.Range("M7").Value = 6.537277325709259
$ws.Range("R7").Value = 9
$ws.Range("S7").Value = 0.04948317572025512

# Row 11 (Indonesian)
$ws.Range("B11").Value = 4477
$ws.Range("C11").Value = 559
$ws.Range("D11").Value = 557
$ws.Range("E11").Value = 30.61000670091579
$ws.Range("F11").Value = 31.83184257602862
$ws.Range("G11").Value = 29.7181328545781
$ws.Range("H11").Value = 4257
$ws.Range("I11").Value = 3.106369626608096
$ws.Range("J11").Value = 2366
$ws.Range("K11").Value = 13.29661683713611
$ws.Range("L11").Value = 2257
$ws.Range("M11").Value = 13.63499063613846
$ws.Range("N11").Value = 168
$ws.Range("O11").Value = 0.122591049393977
$ws.Range("P11").Value = 35
$ws.Range("Q11").Value = 0.1966955153422502
$ws.Range("R11").Value = 21
$ws.Range("S11").Value = 0.1268652208058962

# Row 12 (Finnish)
$ws.Range("D12").Value = 1555
$ws.Range("G12").Value = 29.97041800643087
$ws.Range("L12").Value = 1600
$ws.Range("M12").Value = 3.433181701141533
$ws.Range("R12").Value = 93
$ws.Range("S12").Value = 0.1995536863788516

# Row 15 (Japanese)
$ws.Range("D15").Value = 543
$ws.Range("G15").Value = 35.51565377532228
$ws.Range("L15").Value = 929
$ws.Range("M15").Value = 4.817215452424164
$ws.Range("R15").Value = 63
$ws.Range("S15").Value = 0.3266787658802178
